# "version final sin errores"
# - Bump the Version value (B3) on the Metadata sheet from 0.4.0 to 0.7.0.
# - Remove the "Jurisdiction" / "Chile" property row (row 11) entirely,
#   shifting the rows below it up by one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the Version property value.
$ws1.Cells.Item(3, 2).Value = "0.7.0"

# Delete the whole "Jurisdiction" row (row 11); everything below shifts up.
$ws1.Rows.Item(11).Delete()
